$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> Dll4/Notch2 -> ECs)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 36.75793933333333
$ws.Range("H2").Value = 110.273818
$ws.Range("I2").Value = 0.9858943139827973
$ws.Range("J2").Value = 0.9858943139827971
$ws.Range("M2").Value = 0.9705896666666667
$ws.Range("N2").Value = 2.911769
$ws.Range("O2").Value = 0.02073452941466921
$ws.Range("P2").Value = 0.02073452941466921
$ws.Range("Q2").Value = 35.67687608489356
$ws.Range("R2").Value = 321.091884764042
$ws.Range("S2").Value = 0.02044205465303144
$ws.Range("T2").Value = 0.02044205465303143

# Row 3 (ECs -> Dll4/Notch2 -> FAPs)
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 36.75793933333333
$ws.Range("H3").Value = 110.273818
$ws.Range("I3").Value = 0.9858943139827973
$ws.Range("J3").Value = 0.9858943139827971
$ws.Range("O3").Value = 0.5628689972673966
$ws.Range("P3").Value = 0.5628689972673966
$ws.Range("Q3").Value = 968.5007586104199
$ws.Range("R3").Value = 8716.50682749378
$ws.Range("S3").Value = 0.554929343923125
$ws.Range("T3").Value = 0.5549293439231249

# Row 4 (ECs -> Dll4/Notch2 -> MuSCs)
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 36.75793933333333
$ws.Range("H4").Value = 110.273818
$ws.Range("I4").Value = 0.9858943139827973
$ws.Range("J4").Value = 0.9858943139827971
$ws.Range("M4").Value = 19.49164633333333
$ws.Range("N4").Value = 58.47493899999999
$ws.Range("O4").Value = 0.4163964733179342
$ws.Range("P4").Value = 0.4163964733179341
$ws.Range("Q4").Value = 716.4727534274557
$ws.Range("R4").Value = 6448.254780847102
$ws.Range("S4").Value = 0.4105229154066408
$ws.Range("T4").Value = 0.4105229154066408

# Row 5 (FAPs -> Dll4/Notch2 -> ECs)
$ws.Range("I5").Value = 0.001251989679428792
$ws.Range("J5").Value = 0.001251989679428792
$ws.Range("M5").Value = 0.9705896666666667
$ws.Range("N5").Value = 2.911769
$ws.Range("O5").Value = 0.02073452941466921
$ws.Range("P5").Value = 0.02073452941466921
$ws.Range("Q5").Value = 0.04530615505033334
$ws.Range("R5").Value = 0.407755395453
$ws.Range("S5").Value = 0.00002595941683497858
$ws.Range("T5").Value = 0.00002595941683497857

# Row 6 (FAPs -> Dll4/Notch2 -> FAPs)
$ws.Range("I6").Value = 0.001251989679428792
$ws.Range("J6").Value = 0.001251989679428792
$ws.Range("O6").Value = 0.5628689972673966
$ws.Range("P6").Value = 0.5628689972673966
$ws.Range("S6").Value = 0.0007047061754492137
$ws.Range("T6").Value = 0.0007047061754492134

# Row 7 (FAPs -> Dll4/Notch2 -> MuSCs)
$ws.Range("I7").Value = 0.001251989679428792
$ws.Range("J7").Value = 0.001251989679428792
$ws.Range("M7").Value = 19.49164633333333
$ws.Range("N7").Value = 58.47493899999999
$ws.Range("O7").Value = 0.4163964733179342
$ws.Range("P7").Value = 0.4163964733179341
$ws.Range("Q7").Value = 0.9098505591936665
$ws.Range("R7").Value = 8.188655032742998
$ws.Range("S7").Value = 0.0005213240871446001
$ws.Range("T7").Value = 0.0005213240871445998

# Row 8 (MuSCs -> Dll4/Notch2 -> ECs)
$ws.Range("G8").Value = 0.4792353333333333
$ws.Range("H8").Value = 1.437706
$ws.Range("I8").Value = 0.01285369633777395
$ws.Range("J8").Value = 0.01285369633777395
$ws.Range("M8").Value = 0.9705896666666667
$ws.Range("N8").Value = 2.911769
$ws.Range("O8").Value = 0.02073452941466921
$ws.Range("P8").Value = 0.02073452941466921
$ws.Range("Q8").Value = 0.4651408624348888
$ws.Range("R8").Value = 4.186267761914
$ws.Range("S8").Value = 0.0002665153448028
$ws.Range("T8").Value = 0.0002665153448027999

# Row 9 (MuSCs -> Dll4/Notch2 -> FAPs)
$ws.Range("G9").Value = 0.4792353333333333
$ws.Range("H9").Value = 1.437706
$ws.Range("I9").Value = 0.01285369633777395
$ws.Range("J9").Value = 0.01285369633777395
$ws.Range("O9").Value = 0.5628689972673966
$ws.Range("P9").Value = 0.5628689972673966
$ws.Range("Q9").Value = 12.62692610914
$ws.Range("R9").Value = 113.64233498226
$ws.Range("S9").Value = 0.007234947168822432
$ws.Range("T9").Value = 0.007234947168822431

# Row 10 (MuSCs -> Dll4/Notch2 -> MuSCs)
$ws.Range("G10").Value = 0.4792353333333333
$ws.Range("H10").Value = 1.437706
$ws.Range("I10").Value = 0.01285369633777395
$ws.Range("J10").Value = 0.01285369633777395
$ws.Range("M10").Value = 19.49164633333333
$ws.Range("N10").Value = 58.47493899999999
$ws.Range("O10").Value = 0.4163964733179342
$ws.Range("P10").Value = 0.4163964733179341
$ws.Range("Q10").Value = 9.341085627770443
$ws.Range("R10").Value = 84.06977064993399
$ws.Range("S10").Value = 0.005352233824148719
$ws.Range("T10").Value = 0.005352233824148719
